$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1 - set the text first, then copy the formatting
# from the existing "sum" header (G1) so it matches the other header cells.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New "Save" column data values (0/1 flags) for rows 2-15
$saveValues = @(1, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
